$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 32: consolidate A32:F32 formatting onto the same style already used
# by rows 29-31 (they were duplicate styles in the source file); this frees
# up the now-unused duplicate cell format.
$ws.Range("A29:F29").Copy() | Out-Null
$ws.Range("A32:F32").PasteSpecial(-4122) | Out-Null

# --- New ISFORMULA demo rows (37-40); row 36 stays blank on purpose, it's
# used as the "empty cell" test target below.
$ws.Range("A37").Formula = "=ISFORMULA(A19)"
$ws.Range("B37").Value = "isformula – formula cell"

$ws.Range("A38").Formula = "=ISFORMULA('Foo Bar'!A1)"
$ws.Range("B38").Value = "isformula – number cell"

$ws.Range("A39").Formula = "=ISFORMULA(B37)"
$ws.Range("B39").Value = "isformula – text cell"

$ws.Range("A40").Formula = "=ISFORMULA(A36)"
$ws.Range("B40").Value = "isformula – empty cell"

# Match the number format already used by the neighbouring TRUE()/FALSE() cells.
$ws.Range("A34").Copy() | Out-Null
$ws.Range("A37:A40").PasteSpecial(-4122) | Out-Null

$ws.Range("A41").Select() | Out-Null
